# AOOrderHistory.xlsx -- add separate sheet for TH_TC015
#
# Starting point: a single worksheet "Sheet1" holding test-case rows
# TH_TC014_1..3 for test account test9182.
#
# Target: rename that sheet to TH_TC014, update the 3rd row's total
# price, and add a new worksheet TH_TC015 (placed after TH_TC014) that
# holds a single header + data row for a brand new test account
# (test9183), re-using the same header/value styling as TH_TC014.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------
# 1. Rename the existing sheet
# ---------------------------------------------------------------
$ws1 = $wb.Worksheets.Item(1)
$ws1.Name = "TH_TC014"

# ---------------------------------------------------------------
# 2. Add the new TH_TC015 sheet right after TH_TC014, re-using the
#    header + row styling/values from TH_TC014 (same "TC" /
#    "DT_email" / "DT_password" header, same bold+underline / text
#    formats), then overwrite with the TH_TC015-specific values.
# ---------------------------------------------------------------
$ws2 = $wb.Worksheets.Add($null, $ws1)
$ws2.Name = "TH_TC015"
$ws1.Range("A1:C2").Copy($ws2.Range("A1:C2"))

$ws2.Range("B2").Value = "test9183"

# TH_TC014_3's total price changes from $569.98 to $189.98
$ws1.Range("I4").Value = "$189.98"

$ws2.Range("A2").Value = "TH_TC015_1"
$ws2.Range("C2").Value = "Testing123!"

$ws2.Columns.Item(3).AutoFit() | Out-Null

# ---------------------------------------------------------------
# 3. Selection / active-tab bookkeeping to match the edited file
# ---------------------------------------------------------------
$ws1.Range("C20").Select() | Out-Null
$ws2.Range("H15").Select() | Out-Null
$ws2.Activate() | Out-Null
